# Reposition/resize the two pictures on each of slides 7 and 8
# (sldId 262 and 263) per the "updated the format of my slides" commit.
#
# The point literals below were chosen so that, after PowerPoint's
# internal point->EMU conversion (which narrows to a 32-bit float before
# scaling by 12700 EMU/point and truncating), they reproduce the exact
# target EMU offsets/extents from the authoritative OOXML.

$p = $ppt.ActivePresentation

function Set-ShapeBox($shape, $left, $top, $width, $height) {
    $shape.LockAspectRatio = $false
    $shape.Left = $left
    $shape.Top = $top
    $shape.Width = $width
    $shape.Height = $height
}

# Slide 7 (sldId 262) -----------------------------------------------------
$s7 = $p.Slides.Item(7)

# Picture 3 (cNvPr id="4")
$pic7a = $s7.Shapes.Item(2)
Set-ShapeBox $pic7a 0.0 235.87008666992188 469.0910339355469 301.55853271484375

# Picture 5 (cNvPr id="6")
$pic7b = $s7.Shapes.Item(3)
Set-ShapeBox $pic7b 490.9091491699219 235.87001037597656 469.09088134765625 301.55853271484375

# Slide 8 (sldId 263) -----------------------------------------------------
$s8 = $p.Slides.Item(8)

# Picture 4 (cNvPr id="5")
$pic8a = $s8.Shapes.Item(2)
Set-ShapeBox $pic8a 61.0907096862793 218.4544219970703 418.9093017578125 314.1819763183594

# Picture 6 (cNvPr id="7")
$pic8b = $s8.Shapes.Item(3)
Set-ShapeBox $pic8b 507.272705078125 218.454345703125 418.9093017578125 314.18206787109375
